# Apply odds updates to "Jogos da Semana" worksheet as described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 updates
$ws.Range("K4").Value = 2.3
$ws.Range("N4").Value = 13
$ws.Range("AC4").Value = 12
$ws.Range("AD4").Value = 7.5
$ws.Range("AI4").Value = 15
$ws.Range("AQ4").Value = 41
$ws.Range("AW4").Value = 29

# Row 6 updates
$ws.Range("Q6").Value = 2.1
$ws.Range("R6").Value = 1.7

# Row 7 updates
$ws.Range("G7").Value = 3
$ws.Range("I7").Value = 2.63
$ws.Range("J7").Value = 3.75
$ws.Range("L7").Value = 3.5
$ws.Range("M7").Value = 1.13
$ws.Range("N7").Value = 6
$ws.Range("W7").Value = 7
$ws.Range("Y7").Value = 12
$ws.Range("Z7").Value = 29
$ws.Range("AA7").Value = 29
$ws.Range("AG7").Value = 6.5
$ws.Range("AJ7").Value = 26
$ws.Range("AU7").Value = 4.5
$ws.Range("AV7").Value = 17
$ws.Range("AW7").Value = 34
$ws.Range("BD7").Value = 1250
